$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new historical hammer performance rows (10-12) after the
# existing data (which ends at row 9), following the same column layout:
# B=Performance, C=Date, D=Name, F=Po10 Event, G=Fixture, I=Gender,
# J=Age Code, K=Notes

$ws.Range("B10").Value = "56.42"
$ws.Range("C10").Value = "19 Jun 1993"
$ws.Range("D10").Value = "Nigel Spivey"
$ws.Range("F10").Value = "HT7.26K"
$ws.Range("G10").Value = "Perivale"
$ws.Range("I10").Value = "M"
$ws.Range("J10").Value = "ALL"
$ws.Range("K10").Value = "From Noel Moss 8Apr2024"

$ws.Range("B11").Value = "47.64"
$ws.Range("C11").Value = "11 Sep 1994"
$ws.Range("D11").Value = "Simon Blackwell"
$ws.Range("F11").Value = "HT7.26K"
$ws.Range("G11").Value = "Colchester"
$ws.Range("I11").Value = "M"
$ws.Range("J11").Value = "ALL"
$ws.Range("K11").Value = "From Noel Moss 8Apr2024"

$ws.Range("B12").Value = "47.14"
$ws.Range("C12").Value = "1 Jul 1995"
$ws.Range("D12").Value = "Gary Parsons"
$ws.Range("F12").Value = "HT7.26K"
$ws.Range("G12").Value = "Welwyn"
$ws.Range("I12").Value = "M"
$ws.Range("J12").Value = "ALL"
$ws.Range("K12").Value = "From Noel Moss 8Apr2024"

# Update the view state: scroll so column B is leftmost and select G13
$ws.Activate()
$ws.Range("G13").Select()
